# Applies the "6.0.0" release update to the Insight Result Summary
# StructureDefinition workbook:
#   - Metadata sheet: bump Version, Date, set Publisher, replace the
#     Contact row with a Jurisdiction row, and drop the stray duplicate
#     Contact row beneath it.
#   - Elements sheet: refresh the Short/Definition text for the root
#     Extension row and clear a stale "N/A" RIM mapping cell.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "6.0.0"
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$meta.Range("B9").Value = "Alvearie Team"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 used to duplicate the old "Contact" / "No display for
# ContactDetail" row; remove it so everything below shifts up one row.
$meta.Rows.Item(11).Delete()

$elements = $wb.Worksheets.Item("Elements")

$elements.Range("K2").Value = "Insight Result Summary"
$elements.Range("L2").Value = "Value specific final insight results"
$elements.Range("AJ5").Value = ""
